$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "41.513.15"
Set-TextValue $ws.Range("E2") "  -2.10%  "
Set-TextValue $ws.Range("D3") "2.488.58"
Set-TextValue $ws.Range("E3") "  -0.84%  "
Set-TextValue $ws.Range("E4") "  +0.26%  "
Set-TextValue $ws.Range("D5") "315.06"
Set-TextValue $ws.Range("E5") "  +0.78%  "
Set-TextValue $ws.Range("D6") "94.49"
Set-TextValue $ws.Range("E6") "  -3.82%  "
Set-TextValue $ws.Range("E7") "  -1.80%  "
Set-TextValue $ws.Range("E8") "  +0.19%  "
Set-TextValue $ws.Range("D9") "0.499"
Set-TextValue $ws.Range("E9") "  -2.60%  "
Set-TextValue $ws.Range("D10") "33.61"
Set-TextValue $ws.Range("E10") "  -3.97%  "
Set-TextValue $ws.Range("D11") "0.0784"
Set-TextValue $ws.Range("E11") "  -1.44%  "
Set-TextValue $ws.Range("E12") "  +0.75%  "
Set-TextValue $ws.Range("D13") "2.872.06"
Set-TextValue $ws.Range("E13") "  -0.77%  "
Set-TextValue $ws.Range("D14") "6.93"
Set-TextValue $ws.Range("E14") "  -2.86%  "
Set-TextValue $ws.Range("B15") "WrappedEther"
Set-TextValue $ws.Range("C15") "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D15") "2.552.49"
Set-TextValue $ws.Range("E15") "  +1.92%  "
Set-TextValue $ws.Range("B16") "Chainlink"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D16") "15.46"
Set-TextValue $ws.Range("E16") "  +2.27%  "
Set-TextValue $ws.Range("D17") "0.792"
Set-TextValue $ws.Range("E17") "  -1.14%  "
Set-TextValue $ws.Range("D18") "41.476.50"
Set-TextValue $ws.Range("E18") "  -2.12%  "
Set-TextValue $ws.Range("E19") "  -2.82%  "
Set-TextValue $ws.Range("E20") "  -0.34%  "
Set-TextValue $ws.Range("D21") "70.05"
Set-TextValue $ws.Range("E21") "  +2.55%  "
Set-TextValue $ws.Range("D22") "11.26"
Set-TextValue $ws.Range("E22") "  -5.87%  "
Set-TextValue $ws.Range("D23") "236.86"
Set-TextValue $ws.Range("E23") "  -1.22%  "
Set-TextValue $ws.Range("D24") "2.77"
Set-TextValue $ws.Range("E24") "  -1.92%  "
Set-TextValue $ws.Range("E25") "  -0.03%  "
Set-TextValue $ws.Range("E26") "  -3.62%  "
Set-TextValue $ws.Range("E27") "  -3.90%  "
Set-TextValue $ws.Range("D28") "2.27"
Set-TextValue $ws.Range("E28") "  +1.03%  "
Set-TextValue $ws.Range("E29") "  -0.74%  "
Set-TextValue $ws.Range("D30") "37.18"
Set-TextValue $ws.Range("E30") "  -0.79%  "
Set-TextValue $ws.Range("D31") "154.31"
Set-TextValue $ws.Range("E31") "  -1.29%  "
Set-TextValue $ws.Range("E32") "  -4.49%  "
Set-TextValue $ws.Range("D33") "2.57"
Set-TextValue $ws.Range("E33") "  -2.68%  "
Set-TextValue $ws.Range("D34") "0.0758"
Set-TextValue $ws.Range("E34") "  -2.40%  "
Set-TextValue $ws.Range("D35") "18.02"
Set-TextValue $ws.Range("E35") "  +3.25%  "
Set-TextValue $ws.Range("D36") "3.08"
Set-TextValue $ws.Range("E36") "  -1.05%  "
Set-TextValue $ws.Range("E37") "  -10.49%  "
Set-TextValue $ws.Range("D38") "1.87"
Set-TextValue $ws.Range("E38") "  -3.02%  "
Set-TextValue $ws.Range("E39") "  -1.39%  "
Set-TextValue $ws.Range("E40") "  -4.72%  "
Set-TextValue $ws.Range("E41") "  +0.37%  "
Set-TextValue $ws.Range("E42") "  +0.37%  "
Set-TextValue $ws.Range("D43") "19.85"
Set-TextValue $ws.Range("E43") "  -6.84%  "
Set-TextValue $ws.Range("D44") "1.990.21"
Set-TextValue $ws.Range("E44") "  -0.31%  "
Set-TextValue $ws.Range("E45") "  -1.97%  "
Set-TextValue $ws.Range("E46") "  -5.00%  "
Set-TextValue $ws.Range("D47") "8.85"
Set-TextValue $ws.Range("E47") "  +0.03%  "
Set-TextValue $ws.Range("D48") "2.736.39"
Set-TextValue $ws.Range("E48") "  -0.34%  "
Set-TextValue $ws.Range("D49") "69.66"
Set-TextValue $ws.Range("E49") "  -1.83%  "
Set-TextValue $ws.Range("D50") "97.33"
Set-TextValue $ws.Range("E50") "  -1.83%  "
Set-TextValue $ws.Range("E51") "  -4.09%  "
